$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing last row (row 55, "01-06-2021") with the revised
#     "Actualización desde MV -datos-" figures ---
$ws.Cells.Item(55, 2).Value = 6
$ws.Cells.Item(55, 3).Value = 2.1
$ws.Cells.Item(55, 4).Value = 6.4

# --- Append the new monthly row (row 56, "01-07-2021") ---
# The date label is entered as a quoted-text formula so it is stored as a
# plain string (matching every other "Serie" label in column A) instead of
# being auto-recognised as a date serial number. Copy / PasteSpecial-values
# then collapses the formula down to its literal text result so the saved
# cell is a plain shared-string cell with no formula and no new style.
$ws.Cells.Item(56, 1).Formula = '="01-07-2021"'
$ws.Range("A56").Copy()
$ws.Range("A56").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Cells.Item(56, 2).Value = 6.5
$ws.Cells.Item(56, 3).Value = 1.9
$ws.Cells.Item(56, 4).Value = 6.6
